# Update temp and weather extract raster
#
# Edits applied to the "covars" worksheet (Table1, range A1:G35):
#  1) Row 13 <-> Row 14 swap their data (the "Precipitation" scaled-cluster variable
#     and the "Temperature" cluster variable change places in the sort order),
#     including which of risk_factor_raw (F) / risk_factor_model (G) is flagged "y",
#     and row 14's custom (taller) row height moves along with the swap logic below.
#  2) Row 15 (Temperature, scaled) gains a "y" flag in the risk_factor_model (G) column
#     and becomes the taller (custom height 27) row that row 14 used to be.
#  3) Row 18's column_name changes from "wtrdist_fctb_clst" to
#     "wtrdist_cont_log_scale_clst" (water-distance variable renamed/rescaled).
#  4) The active selection moves to A15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("covars")
$ws.Activate()

# --- Row 13 / Row 14 content swap ---
$ws.Range("A13").Value = "precip_mean_cont_scale_clst"
$ws.Range("B13").Value = "Lagged Precipitation"
$ws.Range("C13").Value = "Precipitation"
$ws.Range("F13").Value = $null
$ws.Range("G13").Value = "y"

$ws.Range("A14").Value = "temp_mean_cont_clst"
$ws.Range("B14").Value = "Lagged Temperature"
$ws.Range("C14").Value = "Temperature"
$ws.Range("F14").Value = "y"
$ws.Range("G14").Value = $null

# Row 14 reverts to the standard (non-custom) row height ...
$ws.Rows.Item(14).AutoFit()

# --- Row 15 gains the "y" flag in risk_factor_model, and the custom row height ---
$ws.Range("G15").Value = "y"
$ws.Rows.Item(15).RowHeight = 27

# --- Row 18 column_name rename ---
$ws.Range("A18").Value = "wtrdist_cont_log_scale_clst"

# --- Selection moves to A15 ---
$ws.Range("A15").Select() | Out-Null
